# Update the Hp-Itgam LR-pairs sheet with newly computed TPM-based values and
# drop the now-obsolete MuSCs -> ECs row (former row 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Hp/Itgam -> MuSCs -------------------------------------
$ws.Range("D2").Value2 = "MuSCs"

$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.1498226666666667
$ws.Range("H2").Value2 = 0.449468
$ws.Range("I2").Value2 = 0.1392821207592237
$ws.Range("J2").Value2 = 0.1392821207592237
$ws.Range("M2").Value2 = 0.01393633333333333
$ws.Range("N2").Value2 = 0.041809
$ws.Range("Q2").Value2 = 0.002087978623555556
$ws.Range("R2").Value2 = 0.018791807612
$ws.Range("S2").Value2 = 0.1392821207592237
$ws.Range("T2").Value2 = 0.1392821207592237

# --- Row 3: FAPs -> Hp/Itgam -> MuSCs ------------------------------------
$ws.Range("D3").Value2 = "MuSCs"

$ws.Range("G3").Value2 = 0.925855
$ws.Range("I3").Value2 = 0.8607178792407764
$ws.Range("J3").Value2 = 0.8607178792407763
$ws.Range("M3").Value2 = 0.01393633333333333
$ws.Range("N3").Value2 = 0.041809
$ws.Range("Q3").Value2 = 0.01290302389833333
$ws.Range("R3").Value2 = 0.116127215085
$ws.Range("S3").Value2 = 0.8607178792407764
$ws.Range("T3").Value2 = 0.8607178792407763

# --- Row 4 (MuSCs -> ECs) is no longer produced by the pipeline: remove it
$ws.Rows.Item(4).Delete()
